$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 16
$ws.Range("C3").Value = 40
$ws.Range("D3").Value = 31
$ws.Range("C4").Value = 81
$ws.Range("D4").Value = 59
$ws.Range("C5").Value = 68
$ws.Range("D5").Value = 62
$ws.Range("C6").Value = 104
$ws.Range("D6").Value = 79
$ws.Range("C7").Value = 69
$ws.Range("D7").Value = 52
$ws.Range("C8").Value = 44
$ws.Range("D8").Value = 29
$ws.Range("C10").Value = 27
$ws.Range("D10").Value = 20
$ws.Range("C11").Value = 60
$ws.Range("D11").Value = 49
$ws.Range("C12").Value = 36
$ws.Range("D12").Value = 28
$ws.Range("C14").Value = 90
$ws.Range("D14").Value = 75
$ws.Range("C15").Value = 54
$ws.Range("D15").Value = 42
$ws.Range("C16").Value = 71
$ws.Range("D16").Value = 56
$ws.Range("C17").Value = 30
$ws.Range("D17").Value = 24
$ws.Range("C18").Value = 43
$ws.Range("D18").Value = 33
$ws.Range("C19").Value = 43
$ws.Range("D19").Value = 33
$ws.Range("C21").Value = 59
$ws.Range("D21").Value = 44
$ws.Range("C22").Value = 31
$ws.Range("D22").Value = 21
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 15
$ws.Range("C24").Value = 63
$ws.Range("D24").Value = 52
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 28
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = 26
$ws.Range("C28").Value = 53
$ws.Range("D28").Value = 39
$ws.Range("C29").Value = 65
$ws.Range("D29").Value = 45
$ws.Range("C30").Value = 50
$ws.Range("D31").Value = 41
$ws.Range("C32").Value = 42
$ws.Range("C33").Value = 47
$ws.Range("D33").Value = 33
$ws.Range("D34").Value = 39
$ws.Range("D35").Value = 41
$ws.Range("C37").Value = 53
$ws.Range("D37").Value = 33
$ws.Range("C38").Value = 62
$ws.Range("D38").Value = 47
$ws.Range("C39").Value = 60
$ws.Range("D39").Value = 51
$ws.Range("D40").Value = 64
$ws.Range("C41").Value = 54
$ws.Range("D41").Value = 42
$ws.Range("D42").Value = 50
$ws.Range("C43").Value = 58
$ws.Range("D43").Value = 45
$ws.Range("C44").Value = 62
$ws.Range("C45").Value = 35
$ws.Range("D45").Value = 20
$ws.Range("C46").Value = 64
$ws.Range("D46").Value = 46
$ws.Range("C47").Value = 47
$ws.Range("D47").Value = 38
$ws.Range("C49").Value = 55
$ws.Range("D49").Value = 45
$ws.Range("C50").Value = 35
$ws.Range("D50").Value = 29
$ws.Range("D51").Value = 28
$ws.Range("C52").Value = 46
$ws.Range("D52").Value = 35
$ws.Range("D53").Value = 69
$ws.Range("C54").Value = 50
$ws.Range("D54").Value = 28
$ws.Range("D55").Value = 42
$ws.Range("D56").Value = 61
$ws.Range("C57").Value = 65
$ws.Range("D57").Value = 39
$ws.Range("C58").Value = 64
$ws.Range("D58").Value = 50
$ws.Range("C59").Value = 60
$ws.Range("D59").Value = 46
$ws.Range("C60").Value = 52
$ws.Range("D60").Value = 32
$ws.Range("C61").Value = 31
$ws.Range("D61").Value = 23
$ws.Range("D62").Value = 38
$ws.Range("C63").Value = 101
$ws.Range("D63").Value = 91
$ws.Range("D64").Value = 22
$ws.Range("C65").Value = 29
$ws.Range("D65").Value = 25
$ws.Range("C66").Value = 60
$ws.Range("D66").Value = 42
$ws.Range("C67").Value = 74
$ws.Range("D67").Value = 58
$ws.Range("C68").Value = 47
$ws.Range("D68").Value = 32
$ws.Range("C69").Value = 94
$ws.Range("D69").Value = 73
$ws.Range("C70").Value = 47
$ws.Range("D71").Value = 59
$ws.Range("C73").Value = 57
$ws.Range("C74").Value = 107
$ws.Range("D74").Value = 91
$ws.Range("D75").Value = 105
$ws.Range("C76").Value = 27
$ws.Range("D76").Value = 14
$ws.Range("C77").Value = 129
$ws.Range("D77").Value = 127
$ws.Range("C78").Value = 74
$ws.Range("D78").Value = 55
$ws.Range("C79").Value = 107
$ws.Range("D79").Value = 105
$ws.Range("C80").Value = 108
$ws.Range("D80").Value = 88
$ws.Range("C81").Value = 72
$ws.Range("D81").Value = 48
$ws.Range("C83").Value = 69
$ws.Range("D83").Value = 59
$ws.Range("C84").Value = 183
$ws.Range("D84").Value = 115
$ws.Range("C92").Value = 219
$ws.Range("D92").Value = 143
$ws.Range("C93").Value = 5110
$ws.Range("D93").Value = 3875
